$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Unprotect()

# Update the confidential disclosure date text in A10
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."

# Update the weight (D) and percent change (E) values for each holding row
$ws.Range("D2").Value = 0.4877218503070156
$ws.Range("E2").Value = -0.002698535080955966

$ws.Range("D3").Value = 0.3311958517023161
$ws.Range("E3").Value = 0.00155369974752384

$ws.Range("D4").Value = 0.09703028593764688
$ws.Range("E4").Value = 0.0005354274495805722

$ws.Range("D5").Value = 0.05401223114488859
$ws.Range("E5").Value = -0.001374098248024813

$ws.Range("D6").Value = 0.03003978090813282
$ws.Range("E6").Value = -0.0009264978381717803

$ws.Range("E7").Value = -0.0008516528373584409

$ws.Protect()
